# Connected Office Test Data - Test Results sheet update
# "Added Read Functionality and published to Orchestrator"
#
# The Read (and related Create/Update/Delete) tests that previously failed
# now pass, so the corresponding boolean cells flip from FALSE to TRUE.
# The stale "bug" note in B26 (and its now-orphaned shared string) is
# removed since the underlying issue has been resolved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Results")
$ws.Activate()

# Rows 2-7, 11-13, 16-21: Update (D) and Delete (E) tests now pass.
$fullRows = @(2,3,4,5,6,7,11,12,13,16,17,18,19,20,21)
foreach ($r in $fullRows) {
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = $true
}

# Rows 8-10, 14-15, 22-24: only the Update (D) test now passes; Delete (E)
# remains failing/unconfirmed.
$updateOnlyRows = @(8,9,10,14,15,22,23,24)
foreach ($r in $updateOnlyRows) {
    $ws.Cells.Item($r, 4).Value = $true
}

# Rows 16-24: Create (B) and Read (C) tests now pass as well.
$createReadRows = @(16,17,18,19,20,21,22,23,24)
foreach ($r in $createReadRows) {
    $ws.Cells.Item($r, 2).Value = $true
    $ws.Cells.Item($r, 3).Value = $true
}

# Remove the stale bug note that lived in row 26 - the update/delete bug
# has been resolved, so clear the cell/row entirely.
$ws.Rows.Item(26).Delete()

# Update the active selection to reflect where the author ended up working.
$ws.Range("B26").Select()
